# "Agregado eventos y funcionalidad varios."
# Insert a new column before column K on the "Resolucion" sheet, give it a
# header "Estado Parada" (new shared string) in row 3 with the same style as
# the other row-3 headers, format the merged-looking row-2 cell above it, and
# tidy up the row heights / selection that Excel recorded after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resolucion")
$ws.Activate()

# Insert a new column K - everything from K onward (incl. merged ranges)
# shifts one column to the right automatically.
$ws.Columns("K").Insert()

# The new column inherits neighbouring column's width category; set it
# explicitly to match column J's width.
$ws.Columns("K").ColumnWidth = 12

# New header cell for the inserted column.
$ws.Range("K3").Value = "Estado Parada"

# Row 2 above the new header: centered, wrapped, no border (matches the
# formatting used for the other row-2 grouping cells, minus the border).
$ws.Range("K2").Borders.LineStyle = -4142
$ws.Range("K2").HorizontalAlignment = -4108
$ws.Range("K2").WrapText = $true

# Row-height touch-ups recorded by Excel after the column insert / edit.
$ws.Rows(3).RowHeight = 75
$ws.Rows(4).RowHeight = 36

# Leave the selection where the user ended up.
$ws.Range("K4").Select()
